$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) ", em que {o/a} " -> ", em que o(a) "  (paragraph 12)
# ------------------------------------------------------------------
$p12 = $d.Paragraphs.Item(12)
$rng = $p12.Range.Duplicate
$rng.Find.Execute("{o/a} ", $true, $false, $false, $false, $false, $true, 1, $false, "o(a) ", 2) | Out-Null

# ------------------------------------------------------------------
# 2) "pertencente ao efetivo {do/da} " -> "pertencente ao efetivo do(a) " (paragraph 12)
# ------------------------------------------------------------------
$p12 = $d.Paragraphs.Item(12)
$rng = $p12.Range.Duplicate
$rng.Find.Execute("{do/da}", $true, $false, $false, $false, $false, $true, 1, $false, "do(a)", 2) | Out-Null

# ------------------------------------------------------------------
# 3) "{O/A} " -> "O(a) " (paragraph 16)
# ------------------------------------------------------------------
$p16 = $d.Paragraphs.Item(16)
$rng = $p16.Range.Duplicate
$rng.Find.Execute("{O/A}", $true, $false, $false, $false, $false, $true, 1, $false, "O(a)", 2) | Out-Null

# ------------------------------------------------------------------
# 4) "{referido/referida} " -> "referido(a) " (paragraph 16)
# ------------------------------------------------------------------
$p16 = $d.Paragraphs.Item(16)
$rng = $p16.Range.Duplicate
$rng.Find.Execute("{referido/referida}", $true, $false, $false, $false, $false, $true, 1, $false, "referido(a)", 2) | Out-Null

# ------------------------------------------------------------------
# 5) "{pronto/pronta} para o " -> "pronto(a) para o " (paragraph 16)
# ------------------------------------------------------------------
$p16 = $d.Paragraphs.Item(16)
$rng = $p16.Range.Duplicate
$rng.Find.Execute("{pronto/pronta} para o ", $true, $false, $false, $false, $false, $true, 1, $false, "pronto(a) para o ", 2) | Out-Null

# ------------------------------------------------------------------
# 6) "À" -> "Ao(à)" only in the "interessada" paragraph (paragraph 18)
# ------------------------------------------------------------------
$p18 = $d.Paragraphs.Item(18)
$rng = $p18.Range.Duplicate
$rng.Find.Execute("À", $true, $false, $false, $false, $false, $true, 1, $false, "Ao(à)", 2) | Out-Null

# ------------------------------------------------------------------
# 7) Move the _GoBack bookmark from the end of paragraph 17 ("... providências;")
#    to the start of paragraph 20 ("{numero_siged}")
# ------------------------------------------------------------------
$p20 = $d.Paragraphs.Item(20)
$bmStart = $p20.Range.Duplicate
$bmStart.Collapse(1) | Out-Null
$d.Bookmarks.Add("_GoBack", $bmStart) | Out-Null

# ------------------------------------------------------------------
# 8) Drop the proofErr wrapping around "numero_siged" while keeping {.../.}
#    as separate runs: merge the whole "{numero_siged}" into a single run,
#    then recolor just the inner text back to EE0000.
# ------------------------------------------------------------------
$p20 = $d.Paragraphs.Item(20)
$rng = $p20.Range.Duplicate
$rng.Find.Execute("{numero_siged}", $true, $false, $false, $false, $false, $true, 1, $false, "{numero_siged}", 2) | Out-Null

$p20 = $d.Paragraphs.Item(20)
$rng2 = $p20.Range.Duplicate
$rng2.Find.Execute("numero_siged") | Out-Null
$rng2.Font.Color = 238
